$wb = $excel.ActiveWorkbook

$wsTestcases = $wb.Worksheets.Item("Testcases")
$wsTime = $wb.Worksheets.Item("Time")

# Fix F7 criteria on Testcases: FAILED -> PASS
$wsTestcases.Range("F7").Value = "PASS"

# Add the missing day of logged time on the Time sheet
$wsTime.Cells.Item(69, 1).Value = 41856
$wsTime.Cells.Item(69, 1).NumberFormat = $wsTime.Cells.Item(68, 1).NumberFormat
$wsTime.Cells.Item(69, 2).Value = "LIBNMATH"
$wsTime.Cells.Item(69, 3).Value = 6

# Update selections / active sheet to match where the edits left the cursor
$wsTestcases.Range("E2").Select()
$wsTime.Activate()
$wsTime.Range("C70").Select()
